$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 22: new match result - Cienciano 3-0 Los Chankas (05/08/2025)
# Force the date to stay as literal text (not auto-converted to a date serial),
# then reset the cell style so no extra number-format style gets attached.
$ws.Range("A22").NumberFormat = "@"
$ws.Range("A22").Value = "05/08/2025"
$ws.Range("A22").Style = "Normal"

$ws.Range("B22").Value = "Cienciano"
$ws.Range("C22").Value = 3
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = "Los Chankas"
$ws.Range("F22").Value = "L"
$ws.Range("G22").Value = 0
$ws.Range("H22").Value = 1
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 2
$ws.Range("K22").Value = 0.23
$ws.Range("L22").Value = 2.76
$ws.Range("M22").Value = 7
$ws.Range("N22").Value = 21
$ws.Range("O22").Value = 2
$ws.Range("P22").Value = 11
